$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.4460406666666667
$ws.Range("H2").Value = 1.338122
$ws.Range("I2").Value = 0.001628842357811546
$ws.Range("J2").Value = 0.001628842357811545
$ws.Range("Q2").Value = 0.02783070739666667
$ws.Range("R2").Value = 0.25047636657
$ws.Range("S2").Value = 0.001628842357811546
$ws.Range("T2").Value = 0.001628842357811545

$ws.Range("G3").Value = 145.8660203333333
$ws.Range("H3").Value = 437.598061
$ws.Range("I3").Value = 0.5326706066061244
$ws.Range("J3").Value = 0.5326706066061244
$ws.Range("Q3").Value = 9.101310338698333
$ws.Range("R3").Value = 81.91179304828499
$ws.Range("S3").Value = 0.5326706066061244
$ws.Range("T3").Value = 0.5326706066061244

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.150912
$ws.Range("H4").Value = 0.452736
$ws.Range("I4").Value = 0.0005510974139175409
$ws.Range("J4").Value = 0.0005510974139175409
$ws.Range("Q4").Value = 0.009416154239999999
$ws.Range("R4").Value = 0.08474538815999999
$ws.Range("S4").Value = 0.0005510974139175409
$ws.Range("T4").Value = 0.0005510974139175409

$ws.Range("G5").Value = 127.376091
$ws.Range("H5").Value = 382.128273
$ws.Range("I5").Value = 0.4651494536221465
$ws.Range("J5").Value = 0.4651494536221465
$ws.Range("Q5").Value = 7.947631197945001
$ws.Range("R5").Value = 71.528680781505
$ws.Range("S5").Value = 0.4651494536221465
$ws.Range("T5").Value = 0.4651494536221465

$wb.Save()
